# Append the 2021 data row (row 12) to Sheet1, mirroring the layout/
# formatting already used by the other yearly rows (e.g. row 11 / 2020年).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$prevRow = 11
$newRow  = 12

# 1. Clone the previous row's formatting (font/border/alignment/number
#    format) onto the new row before writing any values.
$ws.Range("A$prevRow`:U$prevRow").Copy()
$ws.Range("A$newRow`:U$newRow").PasteSpecial(-4122) # xlPasteFormats

# 2. Year label + the numeric series for 2021.
$ws.Range("A$newRow").Value = "2021年"

$values = @{
    "B" = 27805.8
    "D" = 27399.9
    "E" = 67835.60000000001
    "F" = 39415.4
    "G" = 128263.2
    "H" = 154596.3
    "I" = 26390.4
    "J" = 12237.7
    "K" = 26482.3
    "M" = 707379.8
    "N" = 161796.9
    "O" = 82594.8
    "P" = 14013
    "Q" = 23184.7
    "T" = 4694.8
    "U" = 1119629.6
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$newRow").Value = $values[$col]
}

# 3. Columns with no data for this metric stay blank -- like the other
#    rows, they are still present as empty text cells (C, L, R, S).
#    Writing a leading apostrophe forces an empty *text* cell instead of
#    clearing it outright; re-applying the neighbouring row's format
#    afterwards drops the transient quote-prefix style it introduces.
foreach ($col in @("C", "L", "R", "S")) {
    $ws.Range("$col$newRow").Value = "'"
}

$ws.Range("A$prevRow`:U$prevRow").Copy()
$ws.Range("A$newRow`:U$newRow").PasteSpecial(-4122) # xlPasteFormats (re-normalise)

$excel.CutCopyMode = 0
